$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 10
$lastRow = $newRow - 1

# Duplicate the formatting used by the year-label cells (A2:A9) onto the new
# label cell A10 by copying the format from the previous row's label cell.
$ws.Range("A$lastRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10 values (2021年)
$ws.Range("A$newRow").Value = "2021年"
$ws.Range("B$newRow").Value = 3287719
$ws.Range("C$newRow").Value = 1120412
$ws.Range("D$newRow").Value = 7988802
$ws.Range("E$newRow").Value = 21356199
$ws.Range("F$newRow").Value = 774269
$ws.Range("G$newRow").Value = 189665
$ws.Range("H$newRow").Value = 4862351
$ws.Range("I$newRow").Value = 2929454
$ws.Range("J$newRow").Value = 7366985
$ws.Range("K$newRow").Value = 1039905
$ws.Range("L$newRow").Value = 636463
$ws.Range("M$newRow").Value = 61386413
$ws.Range("N$newRow").Value = 239575
$ws.Range("O$newRow").Value = 2552788
$ws.Range("P$newRow").Value = 3329826
$ws.Range("Q$newRow").Value = 2581573
$ws.Range("R$newRow").Value = 1130427
